# Auto-generated Excel COM-interop script to apply numeric updates
# to the Lich_Profits workbook sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(6, 8).Value = 221659.6
$ws.Cells.Item(6, 9).Value = 252074.75
$ws.Cells.Item(6, 10).Value = 99999.0
$ws.Cells.Item(6, 11).Value = 756224.25
$ws.Cells.Item(6, 12).Value = 299997.0
$ws.Cells.Item(6, 13).Value = -756112.25
$ws.Cells.Item(6, 14).Value = -300221.0
$ws.Cells.Item(9, 8).Value = 427.57144
$ws.Cells.Item(9, 9).Value = 434.25
$ws.Cells.Item(9, 11).Value = 434.25
$ws.Cells.Item(9, 13).Value = -265.25
$ws.Cells.Item(28, 8).Value = 1217.5
$ws.Cells.Item(28, 9).Value = 156.25
$ws.Cells.Item(28, 10).Value = 2278.75
$ws.Cells.Item(28, 11).Value = 156.25
$ws.Cells.Item(28, 12).Value = 2278.75
$ws.Cells.Item(28, 13).Value = 328.75
$ws.Cells.Item(28, 14).Value = -3248.75
$ws.Cells.Item(40, 8).Value = 7496.0
$ws.Cells.Item(40, 9).Value = 7496.0
$ws.Cells.Item(40, 10).Value = 0.0
$ws.Cells.Item(40, 11).Value = 7496.0
$ws.Cells.Item(40, 12).Value = 0.0
$ws.Cells.Item(40, 13).Value = -7321.0
$ws.Cells.Item(40, 14).ClearContents() | Out-Null
$ws.Cells.Item(42, 8).Value = 244.83333
$ws.Cells.Item(42, 10).Value = 0.0
$ws.Cells.Item(42, 12).Value = 0.0
$ws.Cells.Item(42, 14).ClearContents() | Out-Null
$ws.Cells.Item(43, 8).Value = 4494.375
$ws.Cells.Item(43, 9).Value = 2279.2856
$ws.Cells.Item(43, 11).Value = 2279.2856
$ws.Cells.Item(43, 13).Value = -2210.2856
$ws.Cells.Item(47, 8).Value = 44750.0
$ws.Cells.Item(47, 9).Value = 44000.0
$ws.Cells.Item(47, 10).Value = 45000.0
$ws.Cells.Item(47, 11).Value = 44000.0
$ws.Cells.Item(47, 12).Value = 45000.0
$ws.Cells.Item(47, 13).Value = -43028.0
$ws.Cells.Item(47, 14).Value = -46944.0
$ws.Cells.Item(53, 8).Value = 497.0
$ws.Cells.Item(53, 10).Value = 609.375
$ws.Cells.Item(53, 12).Value = 609.375
$ws.Cells.Item(53, 14).Value = -1883.375
$ws.Cells.Item(75, 8).Value = 30000.0
$ws.Cells.Item(75, 9).Value = 0.0
$ws.Cells.Item(75, 11).Value = 0.0
$ws.Cells.Item(75, 13).ClearContents() | Out-Null
$ws.Cells.Item(78, 8).Value = 30000.0
$ws.Cells.Item(78, 9).Value = 0.0
$ws.Cells.Item(78, 11).Value = 0.0
$ws.Cells.Item(78, 13).ClearContents() | Out-Null
$ws.Cells.Item(86, 8).Value = 4256.4287
$ws.Cells.Item(86, 9).Value = 4232.5
$ws.Cells.Item(86, 10).Value = 4400.0
$ws.Cells.Item(86, 11).Value = 4232.5
$ws.Cells.Item(86, 12).Value = 4400.0
$ws.Cells.Item(86, 13).Value = -3109.5
$ws.Cells.Item(86, 14).Value = -6646.0
$ws.Cells.Item(89, 8).Value = 4256.4287
$ws.Cells.Item(89, 9).Value = 4232.5
$ws.Cells.Item(89, 10).Value = 4400.0
$ws.Cells.Item(89, 11).Value = 21162.5
$ws.Cells.Item(89, 12).Value = 22000.0
$ws.Cells.Item(89, 13).Value = -15546.5
$ws.Cells.Item(89, 14).Value = -33232.0
$ws.Cells.Item(106, 8).Value = 21199.6
$ws.Cells.Item(106, 9).Value = 1999.6666
$ws.Cells.Item(106, 10).Value = 49999.5
$ws.Cells.Item(106, 11).Value = 1999.6666
$ws.Cells.Item(106, 12).Value = 49999.5
$ws.Cells.Item(106, 13).Value = -1368.6666
$ws.Cells.Item(106, 14).Value = -51261.5
$ws.Cells.Item(113, 8).Value = 8189.5713
$ws.Cells.Item(113, 9).Value = 8894.5
$ws.Cells.Item(113, 11).Value = 8894.5
$ws.Cells.Item(113, 13).Value = -5640.5
$ws.Cells.Item(125, 8).Value = 1303.6
$ws.Cells.Item(125, 9).Value = 931.0
$ws.Cells.Item(125, 11).Value = 8379.0
$ws.Cells.Item(125, 13).Value = -5919.0
$ws.Cells.Item(129, 8).Value = 1226.125
$ws.Cells.Item(129, 9).Value = 1115.7142
$ws.Cells.Item(129, 11).Value = 3347.1426
$ws.Cells.Item(129, 13).Value = 1652.8574
$ws.Cells.Item(132, 8).Value = 7539.2354
$ws.Cells.Item(132, 9).Value = 2866.0
$ws.Cells.Item(132, 10).Value = 18755.0
$ws.Cells.Item(132, 11).Value = 8598.0
$ws.Cells.Item(132, 12).Value = 56265.0
$ws.Cells.Item(132, 13).Value = -6068.0
$ws.Cells.Item(132, 14).Value = -61325.0
$ws.Cells.Item(133, 8).Value = 0.0
$ws.Cells.Item(133, 10).Value = 0.0
$ws.Cells.Item(133, 12).Value = 0.0
$ws.Cells.Item(133, 14).ClearContents() | Out-Null
$ws.Cells.Item(135, 8).Value = 1810.9642
$ws.Cells.Item(135, 9).Value = 1520.6666
$ws.Cells.Item(135, 11).Value = 13685.9994
$ws.Cells.Item(135, 13).Value = -11150.9994
$ws.Cells.Item(138, 8).Value = 3222.0166
$ws.Cells.Item(138, 9).Value = 1750.5333
$ws.Cells.Item(138, 11).Value = 5251.5999
$ws.Cells.Item(138, 13).Value = -111.5999000000002
$ws.Cells.Item(141, 8).Value = 4756.4
$ws.Cells.Item(141, 9).Value = 4756.4
$ws.Cells.Item(141, 11).Value = 14269.2
$ws.Cells.Item(141, 13).Value = -9089.199999999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 7005.742
$ws.Cells.Item(32, 9).Value = 4929.6665
$ws.Cells.Item(32, 10).Value = 20154.223
$ws.Cells.Item(32, 11).Value = 4929.6665
$ws.Cells.Item(32, 12).Value = 20154.223
$ws.Cells.Item(32, 13).Value = -4642.6665
$ws.Cells.Item(32, 14).Value = -20728.223
$ws.Cells.Item(45, 8).Value = 1888.6
$ws.Cells.Item(45, 9).Value = 1926.8572
$ws.Cells.Item(45, 10).Value = 1799.3334
$ws.Cells.Item(45, 11).Value = 1926.8572
$ws.Cells.Item(45, 12).Value = 1799.3334
$ws.Cells.Item(45, 13).Value = -1549.8572
$ws.Cells.Item(45, 14).Value = -2553.3334
$ws.Cells.Item(74, 8).Value = 91704.23
$ws.Cells.Item(74, 9).Value = 117265.12
$ws.Cells.Item(74, 11).Value = 117265.12
$ws.Cells.Item(74, 13).Value = -116391.12
$ws.Cells.Item(77, 8).Value = 91704.23
$ws.Cells.Item(77, 9).Value = 117265.12
$ws.Cells.Item(77, 11).Value = 586325.6
$ws.Cells.Item(77, 13).Value = -581957.6
$ws.Cells.Item(92, 8).Value = 39883.332
$ws.Cells.Item(92, 10).Value = 39883.332
$ws.Cells.Item(92, 12).Value = 39883.332
$ws.Cells.Item(92, 14).Value = -44875.332
$ws.Cells.Item(122, 8).Value = 3236.389
$ws.Cells.Item(122, 9).Value = 2331.1155
$ws.Cells.Item(122, 11).Value = 6993.3465
$ws.Cells.Item(122, 13).Value = -4543.3465
$ws.Cells.Item(132, 8).Value = 2213.1785
$ws.Cells.Item(132, 9).Value = 2183.5186
$ws.Cells.Item(132, 11).Value = 6550.5558
$ws.Cells.Item(132, 13).Value = -4020.5558
$ws.Cells.Item(141, 8).Value = 60000.0
$ws.Cells.Item(141, 10).Value = 60000.0
$ws.Cells.Item(141, 12).Value = 60000.0
$ws.Cells.Item(141, 14).Value = -70360.0

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 5157.5454
$ws.Cells.Item(20, 9).Value = 2403.0
$ws.Cells.Item(20, 10).Value = 8463.0
$ws.Cells.Item(20, 11).Value = 2403.0
$ws.Cells.Item(20, 12).Value = 8463.0
$ws.Cells.Item(20, 13).Value = -2156.0
$ws.Cells.Item(20, 14).Value = -8957.0
$ws.Cells.Item(22, 8).Value = 10528.2
$ws.Cells.Item(22, 9).Value = 381.5
$ws.Cells.Item(22, 11).Value = 381.5
$ws.Cells.Item(22, 13).Value = -208.5
$ws.Cells.Item(82, 8).Value = 100041130.0
$ws.Cells.Item(85, 8).Value = 100041130.0
$ws.Cells.Item(86, 8).Value = 2962.5454
$ws.Cells.Item(86, 9).Value = 2199.0
$ws.Cells.Item(86, 10).Value = 4998.6665
$ws.Cells.Item(86, 11).Value = 2199.0
$ws.Cells.Item(86, 12).Value = 4998.6665
$ws.Cells.Item(86, 13).Value = -1076.0
$ws.Cells.Item(86, 14).Value = -7244.6665
$ws.Cells.Item(89, 8).Value = 2962.5454
$ws.Cells.Item(89, 9).Value = 2199.0
$ws.Cells.Item(89, 10).Value = 4998.6665
$ws.Cells.Item(89, 11).Value = 10995.0
$ws.Cells.Item(89, 12).Value = 24993.3325
$ws.Cells.Item(89, 13).Value = -5379.0
$ws.Cells.Item(89, 14).Value = -36225.3325
$ws.Cells.Item(92, 8).Value = 39945.0
$ws.Cells.Item(92, 10).Value = 39945.0
$ws.Cells.Item(92, 12).Value = 39945.0
$ws.Cells.Item(92, 14).Value = -44937.0
$ws.Cells.Item(96, 8).Value = 250001490.0
$ws.Cells.Item(96, 9).Value = 250001490.0
$ws.Cells.Item(96, 11).Value = 250001490.0
$ws.Cells.Item(96, 13).Value = -249998744.0
$ws.Cells.Item(105, 8).Value = 2112.2104
$ws.Cells.Item(105, 9).Value = 2112.2104
$ws.Cells.Item(105, 11).Value = 2112.2104
$ws.Cells.Item(105, 13).Value = -365.2103999999999

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(4, 8).Value = 250257500.0
$ws.Cells.Item(4, 10).Value = 250257500.0
$ws.Cells.Item(4, 12).Value = 250257500.0
$ws.Cells.Item(4, 14).Value = -250257724.0
$ws.Cells.Item(16, 8).Value = 2671.5454
$ws.Cells.Item(16, 9).Value = 2662.125
$ws.Cells.Item(16, 11).Value = 2662.125
$ws.Cells.Item(16, 13).Value = -2375.125
$ws.Cells.Item(31, 8).Value = 296208.84
$ws.Cells.Item(31, 9).Value = 455800.8
$ws.Cells.Item(31, 10).Value = 3623.5833
$ws.Cells.Item(31, 11).Value = 455800.8
$ws.Cells.Item(31, 12).Value = 3623.5833
$ws.Cells.Item(31, 13).Value = -455505.8
$ws.Cells.Item(31, 14).Value = -4213.5833
$ws.Cells.Item(34, 8).Value = 296208.84
$ws.Cells.Item(34, 9).Value = 455800.8
$ws.Cells.Item(34, 10).Value = 3623.5833
$ws.Cells.Item(34, 11).Value = 455800.8
$ws.Cells.Item(34, 12).Value = 3623.5833
$ws.Cells.Item(34, 13).Value = -455598.8
$ws.Cells.Item(34, 14).Value = -4027.5833
$ws.Cells.Item(58, 8).Value = 2174.5833
$ws.Cells.Item(58, 9).Value = 2097.0557
$ws.Cells.Item(58, 10).Value = 2407.1667
$ws.Cells.Item(58, 11).Value = 2097.0557
$ws.Cells.Item(58, 12).Value = 2407.1667
$ws.Cells.Item(58, 13).Value = -1894.0557
$ws.Cells.Item(58, 14).Value = -2813.1667
$ws.Cells.Item(74, 8).Value = 88578.5
$ws.Cells.Item(74, 10).Value = 88578.5
$ws.Cells.Item(74, 12).Value = 88578.5
$ws.Cells.Item(74, 14).Value = -90326.5
$ws.Cells.Item(77, 8).Value = 88578.5
$ws.Cells.Item(77, 10).Value = 88578.5
$ws.Cells.Item(77, 12).Value = 265735.5
$ws.Cells.Item(77, 14).Value = -274471.5
$ws.Cells.Item(94, 8).Value = 1310.6086
$ws.Cells.Item(94, 9).Value = 746.5
$ws.Cells.Item(94, 11).Value = 746.5
$ws.Cells.Item(94, 13).Value = -295.5
$ws.Cells.Item(99, 8).Value = 363134.0
$ws.Cells.Item(99, 9).Value = 674888.1
$ws.Cells.Item(99, 10).Value = 29111.715
$ws.Cells.Item(99, 11).Value = 674888.1
$ws.Cells.Item(99, 12).Value = 29111.715
$ws.Cells.Item(99, 13).Value = -673390.1
$ws.Cells.Item(99, 14).Value = -32107.715
$ws.Cells.Item(105, 8).Value = 5431.9116
$ws.Cells.Item(105, 9).Value = 2100.4285
$ws.Cells.Item(105, 10).Value = 7763.95
$ws.Cells.Item(105, 11).Value = 2100.4285
$ws.Cells.Item(105, 12).Value = 7763.95
$ws.Cells.Item(105, 13).Value = -353.4285
$ws.Cells.Item(105, 14).Value = -11257.95
$ws.Cells.Item(107, 8).Value = 6474.846
$ws.Cells.Item(107, 9).Value = 1158.6666
$ws.Cells.Item(107, 10).Value = 7168.2607
$ws.Cells.Item(107, 11).Value = 1158.6666
$ws.Cells.Item(107, 12).Value = 7168.2607
$ws.Cells.Item(107, 13).Value = 761.3334
$ws.Cells.Item(107, 14).Value = -11008.2607
$ws.Cells.Item(113, 8).Value = 2671.5454
$ws.Cells.Item(113, 9).Value = 2662.125
$ws.Cells.Item(113, 11).Value = 2662.125
$ws.Cells.Item(113, 13).Value = -492.125
$ws.Cells.Item(122, 8).Value = 3831.0
$ws.Cells.Item(122, 9).Value = 3497.0
$ws.Cells.Item(122, 11).Value = 10491.0
$ws.Cells.Item(122, 13).Value = -8041.0
$ws.Cells.Item(126, 8).Value = 363134.0
$ws.Cells.Item(126, 9).Value = 674888.1
$ws.Cells.Item(126, 10).Value = 29111.715
$ws.Cells.Item(126, 11).Value = 2024664.3
$ws.Cells.Item(126, 12).Value = 87335.145
$ws.Cells.Item(126, 13).Value = -2022194.3
$ws.Cells.Item(126, 14).Value = -92275.145
$ws.Cells.Item(132, 8).Value = 4010.658
$ws.Cells.Item(132, 9).Value = 2036.963
$ws.Cells.Item(132, 10).Value = 8855.182
$ws.Cells.Item(132, 11).Value = 6110.889
$ws.Cells.Item(132, 12).Value = 26565.546
$ws.Cells.Item(132, 13).Value = -3580.889
$ws.Cells.Item(132, 14).Value = -31625.546
$ws.Cells.Item(133, 8).Value = 40504.0
$ws.Cells.Item(133, 10).Value = 40504.0
$ws.Cells.Item(133, 12).Value = 40504.0
$ws.Cells.Item(133, 14).Value = -45564.0
$ws.Cells.Item(134, 8).Value = 5112.46
$ws.Cells.Item(134, 9).Value = 5676.65
$ws.Cells.Item(134, 11).Value = 17029.95
$ws.Cells.Item(134, 13).Value = -14494.95
$ws.Cells.Item(136, 8).Value = 2174.5833
$ws.Cells.Item(136, 9).Value = 2097.0557
$ws.Cells.Item(136, 10).Value = 2407.1667
$ws.Cells.Item(136, 11).Value = 6291.1671
$ws.Cells.Item(136, 12).Value = 7221.500100000001
$ws.Cells.Item(136, 13).Value = -3741.1671
$ws.Cells.Item(136, 14).Value = -12321.5001
$ws.Cells.Item(140, 8).Value = 64888.777
$ws.Cells.Item(140, 10).Value = 68624.875
$ws.Cells.Item(140, 12).Value = 68624.875
$ws.Cells.Item(140, 14).Value = -78984.875

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(5, 8).Value = 2000.6666
$ws.Cells.Item(5, 9).Value = 2004.0
$ws.Cells.Item(5, 10).Value = 1999.0
$ws.Cells.Item(5, 11).Value = 6012.0
$ws.Cells.Item(5, 12).Value = 5997.0
$ws.Cells.Item(5, 13).Value = -5900.0
$ws.Cells.Item(5, 14).Value = -6221.0
$ws.Cells.Item(23, 8).Value = 399.1579
$ws.Cells.Item(23, 10).Value = 484.64285
$ws.Cells.Item(23, 12).Value = 1453.92855
$ws.Cells.Item(23, 14).Value = -1923.92855
$ws.Cells.Item(60, 8).Value = 1128.8
$ws.Cells.Item(60, 9).Value = 1128.8
$ws.Cells.Item(60, 11).Value = 3386.4
$ws.Cells.Item(60, 13).Value = -3135.4
$ws.Cells.Item(116, 8).Value = 14014.5
$ws.Cells.Item(116, 9).Value = 18030.5
$ws.Cells.Item(116, 11).Value = 54091.5
$ws.Cells.Item(116, 13).Value = -50649.5
$ws.Cells.Item(121, 9).Value = 90909870.0
$ws.Cells.Item(121, 10).Value = 2564.4
$ws.Cells.Item(121, 11).Value = 272729610.0
$ws.Cells.Item(121, 12).Value = 7693.200000000001
$ws.Cells.Item(121, 13).Value = -272728300.0
$ws.Cells.Item(121, 14).Value = -10313.2
$ws.Cells.Item(122, 8).Value = 1071.25
$ws.Cells.Item(122, 10).Value = 953.75
$ws.Cells.Item(122, 12).Value = 8583.75
$ws.Cells.Item(122, 14).Value = -13483.75
$ws.Cells.Item(131, 8).Value = 7043894.0
$ws.Cells.Item(131, 10).Value = 1674.1904
$ws.Cells.Item(131, 12).Value = 5022.5712
$ws.Cells.Item(131, 14).Value = -15102.5712
$ws.Cells.Item(135, 8).Value = 2000.6666
$ws.Cells.Item(135, 9).Value = 2004.0
$ws.Cells.Item(135, 10).Value = 1999.0
$ws.Cells.Item(135, 11).Value = 18036.0
$ws.Cells.Item(135, 12).Value = 17991.0
$ws.Cells.Item(135, 13).Value = -15501.0
$ws.Cells.Item(135, 14).Value = -23061.0
$ws.Cells.Item(136, 8).Value = 1962.4242
$ws.Cells.Item(136, 9).Value = 1962.4242
$ws.Cells.Item(136, 11).Value = 5887.2726
$ws.Cells.Item(136, 13).Value = -787.2726000000002
$ws.Cells.Item(139, 8).Value = 2214.7856
$ws.Cells.Item(139, 9).Value = 1532.8
$ws.Cells.Item(139, 10).Value = 3919.75
$ws.Cells.Item(139, 11).Value = 4598.4
$ws.Cells.Item(139, 12).Value = 11759.25
$ws.Cells.Item(139, 13).Value = 541.6000000000004
$ws.Cells.Item(139, 14).Value = -22039.25

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 9388.556
$ws.Cells.Item(70, 9).Value = 9916.5
$ws.Cells.Item(70, 10).Value = 8332.667
$ws.Cells.Item(70, 11).Value = 9916.5
$ws.Cells.Item(70, 12).Value = 8332.667
$ws.Cells.Item(70, 13).Value = -9646.5
$ws.Cells.Item(70, 14).Value = -8872.667
$ws.Cells.Item(73, 8).Value = 9388.556
$ws.Cells.Item(73, 9).Value = 9916.5
$ws.Cells.Item(73, 10).Value = 8332.667
$ws.Cells.Item(73, 11).Value = 9916.5
$ws.Cells.Item(73, 12).Value = 8332.667
$ws.Cells.Item(73, 13).Value = -8980.5
$ws.Cells.Item(73, 14).Value = -10204.667
$ws.Cells.Item(80, 8).Value = 5018.364
$ws.Cells.Item(80, 9).Value = 4595.4
$ws.Cells.Item(80, 10).Value = 5370.8335
$ws.Cells.Item(80, 11).Value = 4595.4
$ws.Cells.Item(80, 12).Value = 5370.8335
$ws.Cells.Item(80, 13).Value = -3597.4
$ws.Cells.Item(80, 14).Value = -7366.8335
$ws.Cells.Item(83, 8).Value = 5018.364
$ws.Cells.Item(83, 9).Value = 4595.4
$ws.Cells.Item(83, 10).Value = 5370.8335
$ws.Cells.Item(83, 11).Value = 22977.0
$ws.Cells.Item(83, 12).Value = 26854.1675
$ws.Cells.Item(83, 13).Value = -17985.0
$ws.Cells.Item(83, 14).Value = -36838.1675
$ws.Cells.Item(97, 8).Value = 4112.5654
$ws.Cells.Item(97, 9).Value = 3460.8
$ws.Cells.Item(97, 11).Value = 3460.8
$ws.Cells.Item(97, 13).Value = -2964.8
$ws.Cells.Item(107, 8).Value = 1235.375
$ws.Cells.Item(107, 9).Value = 1557.6
$ws.Cells.Item(107, 11).Value = 1557.6
$ws.Cells.Item(107, 13).Value = 362.4000000000001
$ws.Cells.Item(113, 8).Value = 36000.0
$ws.Cells.Item(113, 9).Value = 0.0
$ws.Cells.Item(113, 10).Value = 36000.0
$ws.Cells.Item(113, 11).Value = 0.0
$ws.Cells.Item(113, 12).Value = 36000.0
$ws.Cells.Item(113, 13).ClearContents() | Out-Null
$ws.Cells.Item(113, 14).Value = -40340.0
$ws.Cells.Item(122, 8).Value = 4714.4
$ws.Cells.Item(122, 9).Value = 4575.2354
$ws.Cells.Item(122, 11).Value = 13725.7062
$ws.Cells.Item(122, 13).Value = -11275.7062
$ws.Cells.Item(132, 8).Value = 35666.344
$ws.Cells.Item(132, 9).Value = 41370.668
$ws.Cells.Item(132, 10).Value = 4863.0
$ws.Cells.Item(132, 11).Value = 124112.004
$ws.Cells.Item(132, 12).Value = 14589.0
$ws.Cells.Item(132, 13).Value = -121582.004
$ws.Cells.Item(132, 14).Value = -19649.0

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(7, 8).Value = 6568.5
$ws.Cells.Item(7, 9).Value = 6679.7334
$ws.Cells.Item(7, 11).Value = 6679.7334
$ws.Cells.Item(7, 13).Value = -6567.7334
$ws.Cells.Item(16, 8).Value = 16667816.0
$ws.Cells.Item(16, 9).Value = 17858250.0
$ws.Cells.Item(16, 11).Value = 17858250.0
$ws.Cells.Item(16, 13).Value = -17858080.0
$ws.Cells.Item(40, 8).Value = 12864.481
$ws.Cells.Item(40, 9).Value = 13949.913
$ws.Cells.Item(40, 11).Value = 13949.913
$ws.Cells.Item(40, 13).Value = -13813.913
$ws.Cells.Item(46, 8).Value = 6173.909
$ws.Cells.Item(46, 9).Value = 4864.125
$ws.Cells.Item(46, 10).Value = 9666.667
$ws.Cells.Item(46, 11).Value = 4864.125
$ws.Cells.Item(46, 12).Value = 9666.667
$ws.Cells.Item(46, 13).Value = -4676.125
$ws.Cells.Item(46, 14).Value = -10042.667
$ws.Cells.Item(68, 8).Value = 7870.8335
$ws.Cells.Item(68, 10).Value = 2800.0
$ws.Cells.Item(68, 12).Value = 2800.0
$ws.Cells.Item(68, 14).Value = -4298.0
$ws.Cells.Item(71, 8).Value = 7870.8335
$ws.Cells.Item(71, 10).Value = 2800.0
$ws.Cells.Item(71, 12).Value = 14000.0
$ws.Cells.Item(71, 14).Value = -21488.0
$ws.Cells.Item(93, 8).Value = 2001.6522
$ws.Cells.Item(93, 9).Value = 1842.6666
$ws.Cells.Item(93, 10).Value = 2299.75
$ws.Cells.Item(93, 11).Value = 1842.6666
$ws.Cells.Item(93, 12).Value = 2299.75
$ws.Cells.Item(93, 13).Value = -594.6666
$ws.Cells.Item(93, 14).Value = -4795.75
$ws.Cells.Item(100, 8).Value = 6669638.5
$ws.Cells.Item(100, 9).Value = 7695265.0
$ws.Cells.Item(100, 10).Value = 3065.0
$ws.Cells.Item(100, 11).Value = 7695265.0
$ws.Cells.Item(100, 12).Value = 3065.0
$ws.Cells.Item(100, 13).Value = -7694724.0
$ws.Cells.Item(100, 14).Value = -4147.0
$ws.Cells.Item(126, 8).Value = 6568.5
$ws.Cells.Item(126, 9).Value = 6679.7334
$ws.Cells.Item(126, 11).Value = 20039.2002
$ws.Cells.Item(126, 13).Value = -17569.2002
$ws.Cells.Item(132, 8).Value = 6806.39
$ws.Cells.Item(132, 9).Value = 8066.6113
$ws.Cells.Item(132, 11).Value = 24199.8339
$ws.Cells.Item(132, 13).Value = -21669.8339
$ws.Cells.Item(139, 8).Value = 48995.668
$ws.Cells.Item(139, 9).Value = 48999.0
$ws.Cells.Item(139, 10).Value = 48994.0
$ws.Cells.Item(139, 11).Value = 48999.0
$ws.Cells.Item(139, 12).Value = 48994.0
$ws.Cells.Item(139, 13).Value = -43859.0
$ws.Cells.Item(139, 14).Value = -59274.0

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(5, 8).Value = 99999.0
$ws.Cells.Item(5, 10).Value = 99999.0
$ws.Cells.Item(5, 12).Value = 99999.0
$ws.Cells.Item(5, 14).Value = -100223.0
$ws.Cells.Item(11, 8).Value = 4199.5
$ws.Cells.Item(11, 9).Value = 0.0
$ws.Cells.Item(11, 10).Value = 4199.5
$ws.Cells.Item(11, 11).Value = 0.0
$ws.Cells.Item(11, 12).Value = 4199.5
$ws.Cells.Item(11, 13).ClearContents() | Out-Null
$ws.Cells.Item(11, 14).Value = -4483.5
$ws.Cells.Item(57, 8).Value = 59900.0
$ws.Cells.Item(57, 10).Value = 59900.0
$ws.Cells.Item(57, 12).Value = 59900.0
$ws.Cells.Item(57, 14).Value = -61408.0
$ws.Cells.Item(74, 8).Value = 33798.8
$ws.Cells.Item(74, 9).Value = 19000.0
$ws.Cells.Item(74, 10).Value = 37498.5
$ws.Cells.Item(74, 11).Value = 19000.0
$ws.Cells.Item(74, 12).Value = 37498.5
$ws.Cells.Item(74, 13).Value = -18064.0
$ws.Cells.Item(74, 14).Value = -39370.5
$ws.Cells.Item(75, 8).Value = 0.0
$ws.Cells.Item(75, 9).Value = 0.0
$ws.Cells.Item(75, 11).Value = 0.0
$ws.Cells.Item(75, 13).ClearContents() | Out-Null
$ws.Cells.Item(77, 8).Value = 33798.8
$ws.Cells.Item(77, 9).Value = 19000.0
$ws.Cells.Item(77, 10).Value = 37498.5
$ws.Cells.Item(77, 11).Value = 57000.0
$ws.Cells.Item(77, 12).Value = 112495.5
$ws.Cells.Item(77, 13).Value = -52320.0
$ws.Cells.Item(77, 14).Value = -121855.5
$ws.Cells.Item(78, 8).Value = 0.0
$ws.Cells.Item(78, 9).Value = 0.0
$ws.Cells.Item(78, 11).Value = 0.0
$ws.Cells.Item(78, 13).ClearContents() | Out-Null
$ws.Cells.Item(122, 8).Value = 21431.666
$ws.Cells.Item(122, 9).Value = 3144.7
$ws.Cells.Item(122, 10).Value = 58005.6
$ws.Cells.Item(122, 11).Value = 9434.099999999999
$ws.Cells.Item(122, 12).Value = 174016.8
$ws.Cells.Item(122, 13).Value = -6984.099999999999
$ws.Cells.Item(122, 14).Value = -178916.8
$ws.Cells.Item(126, 8).Value = 2553.0
$ws.Cells.Item(126, 9).Value = 2532.5881
$ws.Cells.Item(126, 11).Value = 7597.7643
$ws.Cells.Item(126, 13).Value = -5127.7643
$ws.Cells.Item(132, 8).Value = 1324.9333
$ws.Cells.Item(132, 9).Value = 1072.9166
$ws.Cells.Item(132, 11).Value = 3218.7498
$ws.Cells.Item(132, 13).Value = -688.7498
$ws.Cells.Item(136, 8).Value = 557223.75
$ws.Cells.Item(136, 9).Value = 589860.1
$ws.Cells.Item(136, 11).Value = 1769580.3
$ws.Cells.Item(136, 13).Value = -1767030.3
$ws.Cells.Item(138, 8).Value = 99999.0
$ws.Cells.Item(138, 10).Value = 99999.0
$ws.Cells.Item(138, 12).Value = 99999.0
$ws.Cells.Item(138, 14).Value = -110279.0
